$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.92"
$ws.Range("E2").Value = "'-1.22%"
$ws.Range("E3").Value = "'-2.19%"
$ws.Range("D4").Value = "'4.874"
$ws.Range("E4").Value = "'1.63%"
$ws.Range("E5").Value = "'-0.26%"
$ws.Range("D6").Value = "'6.928"
$ws.Range("E6").Value = "'-0.13%"
$ws.Range("D7").Value = "'1.279"
$ws.Range("E7").Value = "'32.89%"
$ws.Range("D8").Value = "'0.8754"
$ws.Range("E8").Value = "'-0.51%"
$ws.Range("D9").Value = "'0.1556"
$ws.Range("E9").Value = "'5.29%"
$ws.Range("D10").Value = "'0.05105"
$ws.Range("E10").Value = "'-0.60%"
$ws.Range("D11").Value = "'0.07496"
$ws.Range("E11").Value = "'2.48%"
$ws.Range("D12").Value = "'0.02964"
$ws.Range("E12").Value = "'-5.42%"
$ws.Range("D13").Value = "'0.09059"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.001578"
$ws.Range("E14").Value = "'0.81%"
$ws.Range("E15").Value = "'0.96%"
$ws.Range("D16").Value = "'0.005958"
$ws.Range("E16").Value = "'0.84%"
$ws.Range("D17").Value = "'3.454"
$ws.Range("E17").Value = "'0.05%"
$ws.Range("D18").Value = "'3.320"
$ws.Range("E18").Value = "'-2.43%"
$ws.Range("E20").Value = "'0.29%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'1.69%"
$ws.Range("D22").Value = "'3.941"
$ws.Range("E22").Value = "'1.60%"
$ws.Range("D23").Value = "'0.04370"
$ws.Range("E23").Value = "'1.38%"
$ws.Range("D24").Value = "'0.001161"
$ws.Range("E24").Value = "'-1.44%"
$ws.Range("D25").Value = "'0.004210"
$ws.Range("E25").Value = "'-1.92%"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("E27").Value = "'-4.35%"
$ws.Range("D40").Value = "'0.04096"
$ws.Range("E40").Value = "'-0.02%"
$ws.Range("D41").Value = "'0.007026"
$ws.Range("E41").Value = "'3.67%"
$ws.Range("E42").Value = "'0.76%"
$ws.Range("D43").Value = "'0.002192"
$ws.Range("E43").Value = "'-0.37%"
$ws.Range("D44").Value = "'0.01126"
$ws.Range("E44").Value = "'-12.63%"
$ws.Range("D45").Value = "'0.00005223"
$ws.Range("E45").Value = "'0.09%"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "'1.490"
$ws.Range("E46").Value = "'-37.34%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.02000"
$ws.Range("E47").Value = "'-11.12%"
